# Wording changes from PM&C.
# Split the ED performance description cell into two differently-formatted
# runs of rich text, and touch up a couple of cosmetic sheet properties
# that Excel recomputes as a side effect of the edit.

$wb = $excel.ActiveWorkbook

$wsDescription = $wb.Worksheets.Item("Description")
$wsData = $wb.Worksheets.Item("Data")

$run1 = 'Emergency department (ED) performance is based on the percentage of people who are seen within benchmark times for five triage categories—triage category 1 is for those with life threatening conditions to be seen "within seconds" and triage category 5 is for those with the least urgent conditions to be seen "within 120 minutes". '
$run2 = 'Performance in triage categories 4 and 5 is likely to be affected by the availability of primary care, which is largely a Commonwealth Government responsibility.'

$cell = $wsDescription.Range("B5")
$cell.Value = $run1 + $run2

$len1 = $run1.Length
$len2 = $run2.Length

$chars1 = $cell.Characters(1, $len1)
$chars1.Font.Name = "Arial"
$chars1.Font.Size = 10
$chars1.Font.Color = 0

$chars2 = $cell.Characters($len1 + 1, $len2)
$chars2.Font.Name = "Calibri"
$chars2.Font.Size = 12
$chars2.Font.Color = 0

# Row 5 grows taller to fit the extra sentence.
$wsDescription.Rows.Item(5).RowHeight = 59.2

# Selection left where the editor's cursor landed after the edit.
$wsDescription.Range("B5").Select()

# Minor column-width rounding nudge on the Data sheet (side effect of the
# font substitutions made above triggering a width recalculation).
$wsData.Columns.Item(1).ColumnWidth = 7.7
